$d = $word.ActiveDocument

$d.Content.Find.Execute("25+52=77", $true, $false, $false, $false, $false, $true, 1, $false, "30+39=69", 2) | Out-Null
$d.Content.Find.Execute("26+54=80", $true, $false, $false, $false, $false, $true, 1, $false, "20+18=38", 2) | Out-Null
$d.Content.Find.Execute("99-80=19", $true, $false, $false, $false, $false, $true, 1, $false, "85-51=34", 2) | Out-Null
$d.Content.Find.Execute("0+35=35", $true, $false, $false, $false, $false, $true, 1, $false, "23+75=98", 2) | Out-Null
$d.Content.Find.Execute("34-25=9", $true, $false, $false, $false, $false, $true, 1, $false, "1+80=81", 2) | Out-Null
$d.Content.Find.Execute("5+16=21", $true, $false, $false, $false, $false, $true, 1, $false, "60+3=63", 2) | Out-Null
$d.Content.Find.Execute("24+1=25", $true, $false, $false, $false, $false, $true, 1, $false, "35+56=91", 2) | Out-Null
$d.Content.Find.Execute("24-23=1", $true, $false, $false, $false, $false, $true, 1, $false, "39+0=39", 2) | Out-Null
$d.Content.Find.Execute("1+19=20", $true, $false, $false, $false, $false, $true, 1, $false, "40+9=49", 2) | Out-Null
$d.Content.Find.Execute("0+55=55", $true, $false, $false, $false, $false, $true, 1, $false, "78-7=71", 2) | Out-Null
$d.Content.Find.Execute("55-38=17", $true, $false, $false, $false, $false, $true, 1, $false, "25-12=13", 2) | Out-Null
$d.Content.Find.Execute("45+14=59", $true, $false, $false, $false, $false, $true, 1, $false, "48+13=61", 2) | Out-Null
$d.Content.Find.Execute("66-15=51", $true, $false, $false, $false, $false, $true, 1, $false, "51+23=74", 2) | Out-Null
$d.Content.Find.Execute("93-29=64", $true, $false, $false, $false, $false, $true, 1, $false, "94-3=91", 2) | Out-Null
$d.Content.Find.Execute("66+15=81", $true, $false, $false, $false, $false, $true, 1, $false, "61+35=96", 2) | Out-Null
$d.Content.Find.Execute("0+59=59", $true, $false, $false, $false, $false, $true, 1, $false, "23+65=88", 2) | Out-Null
$d.Content.Find.Execute("26+8=34", $true, $false, $false, $false, $false, $true, 1, $false, "38-37=1", 2) | Out-Null
$d.Content.Find.Execute("4+90=94", $true, $false, $false, $false, $false, $true, 1, $false, "34+32=66", 2) | Out-Null
$d.Content.Find.Execute("22+38=60", $true, $false, $false, $false, $false, $true, 1, $false, "64-34=30", 2) | Out-Null
$d.Content.Find.Execute("95+2=97", $true, $false, $false, $false, $false, $true, 1, $false, "20+23=43", 2) | Out-Null
$d.Content.Find.Execute("16+39=55", $true, $false, $false, $false, $false, $true, 1, $false, "21+50=71", 2) | Out-Null
$d.Content.Find.Execute("43-39=4", $true, $false, $false, $false, $false, $true, 1, $false, "86-56=30", 2) | Out-Null
$d.Content.Find.Execute("20+53=73", $true, $false, $false, $false, $false, $true, 1, $false, "69-26=43", 2) | Out-Null
$d.Content.Find.Execute("4+4=8", $true, $false, $false, $false, $false, $true, 1, $false, "79-47=32", 2) | Out-Null
$d.Content.Find.Execute("92-37=55", $true, $false, $false, $false, $false, $true, 1, $false, "33+42=75", 2) | Out-Null
$d.Content.Find.Execute("63-45=18", $true, $false, $false, $false, $false, $true, 1, $false, "44-22=22", 2) | Out-Null
$d.Content.Find.Execute("29+28=57", $true, $false, $false, $false, $false, $true, 1, $false, "96-54=42", 2) | Out-Null
$d.Content.Find.Execute("30-1=29", $true, $false, $false, $false, $false, $true, 1, $false, "39-18=21", 2) | Out-Null
$d.Content.Find.Execute("50-28=22", $true, $false, $false, $false, $false, $true, 1, $false, "87-29=58", 2) | Out-Null
$d.Content.Find.Execute("18+73=91", $true, $false, $false, $false, $false, $true, 1, $false, "78-54=24", 2) | Out-Null
$d.Content.Find.Execute("89-7=82", $true, $false, $false, $false, $false, $true, 1, $false, "62-24=38", 2) | Out-Null
$d.Content.Find.Execute("60+15=75", $true, $false, $false, $false, $false, $true, 1, $false, "6-4=2", 2) | Out-Null
$d.Content.Find.Execute("3+73=76", $true, $false, $false, $false, $false, $true, 1, $false, "86+11=97", 2) | Out-Null
$d.Content.Find.Execute("33+58=91", $true, $false, $false, $false, $false, $true, 1, $false, "79-27=52", 2) | Out-Null
$d.Content.Find.Execute("46-23=23", $true, $false, $false, $false, $false, $true, 1, $false, "99-77=22", 2) | Out-Null
$d.Content.Find.Execute("3+80=83", $true, $false, $false, $false, $false, $true, 1, $false, "29+27=56", 2) | Out-Null
$d.Content.Find.Execute("42+35=77", $true, $false, $false, $false, $false, $true, 1, $false, "60-24=36", 2) | Out-Null
$d.Content.Find.Execute("10+34=44", $true, $false, $false, $false, $false, $true, 1, $false, "45-2=43", 2) | Out-Null
$d.Content.Find.Execute("80-35=45", $true, $false, $false, $false, $false, $true, 1, $false, "90-4=86", 2) | Out-Null
$d.Content.Find.Execute("53+17=70", $true, $false, $false, $false, $false, $true, 1, $false, "47-39=8", 2) | Out-Null
$d.Content.Find.Execute("66-24=42", $true, $false, $false, $false, $false, $true, 1, $false, "92-51=41", 2) | Out-Null
$d.Content.Find.Execute("4+86=90", $true, $false, $false, $false, $false, $true, 1, $false, "66-61=5", 2) | Out-Null
$d.Content.Find.Execute("90+5=95", $true, $false, $false, $false, $false, $true, 1, $false, "53+44=97", 2) | Out-Null
$d.Content.Find.Execute("28+59=87", $true, $false, $false, $false, $false, $true, 1, $false, "31+38=69", 2) | Out-Null
$d.Content.Find.Execute("91-48=43", $true, $false, $false, $false, $false, $true, 1, $false, "77-64=13", 2) | Out-Null
$d.Content.Find.Execute("7+79=86", $true, $false, $false, $false, $false, $true, 1, $false, "91+8=99", 2) | Out-Null
$d.Content.Find.Execute("57+42=99", $true, $false, $false, $false, $false, $true, 1, $false, "8+38=46", 2) | Out-Null
$d.Content.Find.Execute("58-39=19", $true, $false, $false, $false, $false, $true, 1, $false, "70+10=80", 2) | Out-Null
$d.Content.Find.Execute("61+1=62", $true, $false, $false, $false, $false, $true, 1, $false, "68+7=75", 2) | Out-Null
$d.Content.Find.Execute("3+25=28", $true, $false, $false, $false, $false, $true, 1, $false, "53+32=85", 2) | Out-Null
$d.Content.Find.Execute("8+17=25", $true, $false, $false, $false, $false, $true, 1, $false, "56-2=54", 2) | Out-Null
$d.Content.Find.Execute("15+13=28", $true, $false, $false, $false, $false, $true, 1, $false, "79-1=78", 2) | Out-Null
$d.Content.Find.Execute("52+18=70", $true, $false, $false, $false, $false, $true, 1, $false, "61-43=18", 2) | Out-Null
$d.Content.Find.Execute("94+1=95", $true, $false, $false, $false, $false, $true, 1, $false, "26+3=29", 2) | Out-Null
$d.Content.Find.Execute("27+42=69", $true, $false, $false, $false, $false, $true, 1, $false, "5+54=59", 2) | Out-Null
$d.Content.Find.Execute("41+32=73", $true, $false, $false, $false, $false, $true, 1, $false, "22+25=47", 2) | Out-Null
$d.Content.Find.Execute("69-59=10", $true, $false, $false, $false, $false, $true, 1, $false, "12+40=52", 2) | Out-Null
$d.Content.Find.Execute("64-22=42", $true, $false, $false, $false, $false, $true, 1, $false, "11+51=62", 2) | Out-Null
$d.Content.Find.Execute("1+75=76", $true, $false, $false, $false, $false, $true, 1, $false, "26+64=90", 2) | Out-Null
$d.Content.Find.Execute("72-53=19", $true, $false, $false, $false, $false, $true, 1, $false, "69+4=73", 2) | Out-Null
$d.Content.Find.Execute("27+55=82", $true, $false, $false, $false, $false, $true, 1, $false, "52+9=61", 2) | Out-Null
$d.Content.Find.Execute("79-0=79", $true, $false, $false, $false, $false, $true, 1, $false, "65-23=42", 2) | Out-Null
$d.Content.Find.Execute("71+3=74", $true, $false, $false, $false, $false, $true, 1, $false, "90-60=30", 2) | Out-Null
$d.Content.Find.Execute("23+4=27", $true, $false, $false, $false, $false, $true, 1, $false, "58-1=57", 2) | Out-Null
$d.Content.Find.Execute("60-44=16", $true, $false, $false, $false, $false, $true, 1, $false, "59+9=68", 2) | Out-Null
$d.Content.Find.Execute("3+4=7", $true, $false, $false, $false, $false, $true, 1, $false, "66-46=20", 2) | Out-Null
$d.Content.Find.Execute("46+8=54", $true, $false, $false, $false, $false, $true, 1, $false, "29+23=52", 2) | Out-Null
$d.Content.Find.Execute("69-61=8", $true, $false, $false, $false, $false, $true, 1, $false, "42-28=14", 2) | Out-Null
$d.Content.Find.Execute("16+48=64", $true, $false, $false, $false, $false, $true, 1, $false, "75-28=47", 2) | Out-Null
$d.Content.Find.Execute("76-14=62", $true, $false, $false, $false, $false, $true, 1, $false, "21+9=30", 2) | Out-Null
$d.Content.Find.Execute("87-35=52", $true, $false, $false, $false, $false, $true, 1, $false, "67-52=15", 2) | Out-Null
$d.Content.Find.Execute("93+1=94", $true, $false, $false, $false, $false, $true, 1, $false, "77+16=93", 2) | Out-Null
$d.Content.Find.Execute("60+39=99", $true, $false, $false, $false, $false, $true, 1, $false, "98-79=19", 2) | Out-Null
$d.Content.Find.Execute("71-2=69", $true, $false, $false, $false, $false, $true, 1, $false, "76-10=66", 2) | Out-Null
$d.Content.Find.Execute("2+29=31", $true, $false, $false, $false, $false, $true, 1, $false, "16+76=92", 2) | Out-Null
$d.Content.Find.Execute("26+65=91", $true, $false, $false, $false, $false, $true, 1, $false, "18+53=71", 2) | Out-Null
$d.Content.Find.Execute("90-90=0", $true, $false, $false, $false, $false, $true, 1, $false, "99-4=95", 2) | Out-Null
$d.Content.Find.Execute("62+4=66", $true, $false, $false, $false, $false, $true, 1, $false, "61+12=73", 2) | Out-Null
$d.Content.Find.Execute("30-15=15", $true, $false, $false, $false, $false, $true, 1, $false, "33+5=38", 2) | Out-Null
$d.Content.Find.Execute("87-30=57", $true, $false, $false, $false, $false, $true, 1, $false, "1+82=83", 2) | Out-Null
$d.Content.Find.Execute("20+15=35", $true, $false, $false, $false, $false, $true, 1, $false, "96-0=96", 2) | Out-Null
$d.Content.Find.Execute("33+13=46", $true, $false, $false, $false, $false, $true, 1, $false, "79-2=77", 2) | Out-Null
$d.Content.Find.Execute("0+79=79", $true, $false, $false, $false, $false, $true, 1, $false, "83-25=58", 2) | Out-Null
$d.Content.Find.Execute("68-66=2", $true, $false, $false, $false, $false, $true, 1, $false, "35+54=89", 2) | Out-Null
$d.Content.Find.Execute("82-73=9", $true, $false, $false, $false, $false, $true, 1, $false, "94-66=28", 2) | Out-Null
$d.Content.Find.Execute("47+41=88", $true, $false, $false, $false, $false, $true, 1, $false, "42-9=33", 2) | Out-Null
$d.Content.Find.Execute("81-2=79", $true, $false, $false, $false, $false, $true, 1, $false, "2+15=17", 2) | Out-Null
$d.Content.Find.Execute("39+1=40", $true, $false, $false, $false, $false, $true, 1, $false, "85-71=14", 2) | Out-Null
$d.Content.Find.Execute("27+64=91", $true, $false, $false, $false, $false, $true, 1, $false, "65-11=54", 2) | Out-Null
$d.Content.Find.Execute("37+15=52", $true, $false, $false, $false, $false, $true, 1, $false, "40+17=57", 2) | Out-Null
$d.Content.Find.Execute("48+11=59", $true, $false, $false, $false, $false, $true, 1, $false, "73+23=96", 2) | Out-Null
$d.Content.Find.Execute("20+25=45", $true, $false, $false, $false, $false, $true, 1, $false, "6+35=41", 2) | Out-Null
$d.Content.Find.Execute("97-81=16", $true, $false, $false, $false, $false, $true, 1, $false, "38+8=46", 2) | Out-Null
$d.Content.Find.Execute("42-18=24", $true, $false, $false, $false, $false, $true, 1, $false, "80-9=71", 2) | Out-Null
$d.Content.Find.Execute("7+78=85", $true, $false, $false, $false, $false, $true, 1, $false, "8+54=62", 2) | Out-Null
$d.Content.Find.Execute("62-23=39", $true, $false, $false, $false, $false, $true, 1, $false, "5+71=76", 2) | Out-Null
$d.Content.Find.Execute("7+20=27", $true, $false, $false, $false, $false, $true, 1, $false, "59+26=85", 2) | Out-Null
$d.Content.Find.Execute("40+28=68", $true, $false, $false, $false, $false, $true, 1, $false, "1+89=90", 2) | Out-Null
$d.Content.Find.Execute("9+86=95", $true, $false, $false, $false, $false, $true, 1, $false, "81-32=49", 2) | Out-Null
$d.Content.Find.Execute("97-97=0", $true, $false, $false, $false, $false, $true, 1, $false, "55+44=99", 2) | Out-Null
